$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "hdhewgeyuf"

$ws.Cells.Item(1, 10).Value = 33.58965063095093
$ws.Cells.Item(2, 2).Value = 1870
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 8).Value = 99.73247726056715
$ws.Cells.Item(2, 9).Value = 0.002680965147453083
$ws.Cells.Item(2, 10).Value = 44.35858583450317
$ws.Cells.Item(3, 10).Value = 40.18504548072815
$ws.Cells.Item(4, 2).Value = 2602
$ws.Cells.Item(4, 4).Value = 2567
$ws.Cells.Item(4, 5).Value = 34
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 99.84441851419682
$ws.Cells.Item(4, 8).Value = 98.69281045751634
$ws.Cells.Item(4, 9).Value = 0.01477449455676516
$ws.Cells.Item(4, 10).Value = 34.38748216629028
$ws.Cells.Item(5, 2).Value = 2029
$ws.Cells.Item(5, 4).Value = 2024
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(5, 7).Value = 99.90128331688055
$ws.Cells.Item(5, 8).Value = 99.80276134122288
$ws.Cells.Item(5, 9).Value = 0.002960039467192896
$ws.Cells.Item(5, 10).Value = 41.0810980796814
$ws.Cells.Item(6, 2).Value = 1768
$ws.Cells.Item(6, 4).Value = 1751
$ws.Cells.Item(6, 5).Value = 16
$ws.Cells.Item(6, 6).Value = 11
$ws.Cells.Item(6, 7).Value = 99.37570942111238
$ws.Cells.Item(6, 8).Value = 99.09451046972269
$ws.Cells.Item(6, 9).Value = 0.01531480431083381
$ws.Cells.Item(6, 10).Value = 36.03067946434021
$ws.Cells.Item(7, 10).Value = 39.0495445728302
$ws.Cells.Item(8, 10).Value = 34.26764059066772
$ws.Cells.Item(9, 2).Value = 2540
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 8).Value = 99.96061441512407
$ws.Cells.Item(9, 9).Value = 0.0003938558487593541
$ws.Cells.Item(9, 10).Value = 38.16045951843262
$ws.Cells.Item(10, 2).Value = 1827
$ws.Cells.Item(10, 5).Value = 33
$ws.Cells.Item(10, 8).Value = 98.19277108433735
$ws.Cells.Item(10, 9).Value = 0.01894150417827298
$ws.Cells.Item(10, 10).Value = 41.08096623420715
$ws.Cells.Item(11, 2).Value = 1883
$ws.Cells.Item(11, 4).Value = 1877
$ws.Cells.Item(11, 5).Value = 5
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 99.94675186368477
$ws.Cells.Item(11, 8).Value = 99.73432518597237
$ws.Cells.Item(11, 9).Value = 0.003193187865886109
$ws.Cells.Item(11, 10).Value = 29.68592405319214
$ws.Cells.Item(12, 10).Value = 33.80615139007568
$ws.Cells.Item(13, 2).Value = 2390
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 8).Value = 99.95814148179154
$ws.Cells.Item(13, 9).Value = 0.009950248756218905
$ws.Cells.Item(13, 10).Value = 33.20887303352356
$ws.Cells.Item(14, 2).Value = 1536
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 8).Value = 99.93485342019544
$ws.Cells.Item(14, 9).Value = 0.0006514657980456026
$ws.Cells.Item(14, 10).Value = 32.87627124786377
$ws.Cells.Item(15, 2).Value = 2286
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 8).Value = 99.64989059080963
$ws.Cells.Item(15, 9).Value = 0.003511852502194908
$ws.Cells.Item(15, 10).Value = 34.20742869377136
$ws.Cells.Item(16, 2).Value = 1989
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 8).Value = 99.89939637826961
$ws.Cells.Item(16, 9).Value = 0.001006542526421741
$ws.Cells.Item(16, 10).Value = 36.54154515266418
$ws.Cells.Item(17, 10).Value = 37.71901297569275
$ws.Cells.Item(18, 10).Value = 36.83397769927979
$ws.Cells.Item(19, 2).Value = 1519
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 8).Value = 99.93412384716733
$ws.Cells.Item(19, 9).Value = 0.0006587615283267457
$ws.Cells.Item(19, 10).Value = 36.73185324668884
$ws.Cells.Item(20, 2).Value = 1614
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 8).Value = 99.93800371977682
$ws.Cells.Item(20, 9).Value = 0.004323656578134651
$ws.Cells.Item(20, 10).Value = 28.79729151725769
$ws.Cells.Item(21, 2).Value = 2633
$ws.Cells.Item(21, 4).Value = 2598
$ws.Cells.Item(21, 5).Value = 34
$ws.Cells.Item(21, 6).Value = 2
$ws.Cells.Item(21, 7).Value = 99.92307692307692
$ws.Cells.Item(21, 8).Value = 98.7082066869301
$ws.Cells.Item(21, 9).Value = 0.01384083044982699
$ws.Cells.Item(21, 10).Value = 34.93981552124023
$ws.Cells.Item(22, 2).Value = 1945
$ws.Cells.Item(22, 4).Value = 1944
$ws.Cells.Item(22, 6).Value = 18
$ws.Cells.Item(22, 7).Value = 99.08256880733946
$ws.Cells.Item(22, 9).Value = 0.009169638308711156
$ws.Cells.Item(22, 10).Value = 40.83864545822144
$ws.Cells.Item(23, 2).Value = 2136
$ws.Cells.Item(23, 4).Value = 2134
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 99.95316159250585
$ws.Cells.Item(23, 8).Value = 99.95316159250585
$ws.Cells.Item(23, 9).Value = 0.0009363295880149813
$ws.Cells.Item(23, 10).Value = 31.30275845527649
$ws.Cells.Item(24, 2).Value = 3017
$ws.Cells.Item(24, 4).Value = 2966
$ws.Cells.Item(24, 5).Value = 50
$ws.Cells.Item(24, 6).Value = 13
$ws.Cells.Item(24, 7).Value = 99.5636119503189
$ws.Cells.Item(24, 8).Value = 98.342175066313
$ws.Cells.Item(24, 9).Value = 0.02114093959731544
$ws.Cells.Item(24, 10).Value = 36.60566306114197
$ws.Cells.Item(25, 10).Value = 35.75739860534668
$ws.Cells.Item(26, 2).Value = 1850
$ws.Cells.Item(26, 4).Value = 1837
$ws.Cells.Item(26, 5).Value = 12
$ws.Cells.Item(26, 6).Value = 22
$ws.Cells.Item(26, 7).Value = 98.81656804733728
$ws.Cells.Item(26, 8).Value = 99.35100054083289
$ws.Cells.Item(26, 9).Value = 0.01827956989247312
$ws.Cells.Item(26, 10).Value = 37.48452830314636
$ws.Cells.Item(27, 4).Value = 2940
$ws.Cells.Item(27, 5).Value = 6
$ws.Cells.Item(27, 6).Value = 14
$ws.Cells.Item(27, 7).Value = 99.52606635071091
$ws.Cells.Item(27, 8).Value = 99.79633401221996
$ws.Cells.Item(27, 9).Value = 0.00676818950930626
$ws.Cells.Item(27, 10).Value = 35.51144623756409
$ws.Cells.Item(28, 10).Value = 34.62693881988525
$ws.Cells.Item(29, 2).Value = 2663
$ws.Cells.Item(29, 4).Value = 2644
$ws.Cells.Item(29, 5).Value = 18
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(29, 7).Value = 99.81124952812382
$ws.Cells.Item(29, 8).Value = 99.32381667918858
$ws.Cells.Item(29, 9).Value = 0.008679245283018867
$ws.Cells.Item(29, 10).Value = 37.65887832641602
$ws.Cells.Item(30, 10).Value = 35.753338098526
$ws.Cells.Item(31, 2).Value = 3250
$ws.Cells.Item(31, 4).Value = 3249
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 99.96923076923076
$ws.Cells.Item(31, 9).Value = 0.0003075976622577669
$ws.Cells.Item(31, 10).Value = 36.14133620262146
$ws.Cells.Item(32, 2).Value = 2266
$ws.Cells.Item(32, 4).Value = 2258
$ws.Cells.Item(32, 6).Value = 3
$ws.Cells.Item(32, 7).Value = 99.86731534719151
$ws.Cells.Item(32, 8).Value = 99.69094922737307
$ws.Cells.Item(32, 9).Value = 0.004420866489832007
$ws.Cells.Item(32, 10).Value = 37.47408437728882
$ws.Cells.Item(33, 10).Value = 35.22914171218872
$ws.Cells.Item(34, 10).Value = 38.23012804985046
$ws.Cells.Item(35, 10).Value = 42.97950673103333
$ws.Cells.Item(36, 2).Value = 2428
$ws.Cells.Item(36, 4).Value = 2425
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 99.95877988458368
$ws.Cells.Item(36, 8).Value = 99.91759373712402
$ws.Cells.Item(36, 9).Value = 0.001236093943139679
$ws.Cells.Item(36, 10).Value = 35.18246817588806
$ws.Cells.Item(37, 2).Value = 2487
$ws.Cells.Item(37, 4).Value = 2482
$ws.Cells.Item(37, 5).Value = 4
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 100
$ws.Cells.Item(37, 8).Value = 99.8390989541432
$ws.Cells.Item(37, 9).Value = 0.001610954490535642
$ws.Cells.Item(37, 10).Value = 35.42835903167725
$ws.Cells.Item(38, 10).Value = 30.38698148727417
$ws.Cells.Item(39, 2).Value = 2089
$ws.Cells.Item(39, 4).Value = 2048
$ws.Cells.Item(39, 5).Value = 40
$ws.Cells.Item(39, 6).Value = 4
$ws.Cells.Item(39, 7).Value = 99.80506822612085
$ws.Cells.Item(39, 8).Value = 98.08429118773947
$ws.Cells.Item(39, 9).Value = 0.02143205065757428
$ws.Cells.Item(39, 10).Value = 39.47330570220947
$ws.Cells.Item(40, 2).Value = 2258
$ws.Cells.Item(40, 5).Value = 2
$ws.Cells.Item(40, 8).Value = 99.9113867966327
$ws.Cells.Item(40, 9).Value = 0.0008865248226950354
$ws.Cells.Item(40, 10).Value = 40.66242289543152
$ws.Cells.Item(41, 10).Value = 36.87566328048706
$ws.Cells.Item(42, 2).Value = 1782
$ws.Cells.Item(42, 4).Value = 1779
$ws.Cells.Item(42, 5).Value = 2
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 100
$ws.Cells.Item(42, 8).Value = 99.88770353733858
$ws.Cells.Item(42, 9).Value = 0.001123595505617978
$ws.Cells.Item(42, 10).Value = 36.01734900474548
$ws.Cells.Item(43, 10).Value = 40.97216868400574
$ws.Cells.Item(44, 10).Value = 37.22258543968201
